$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 250, shifting existing rows 250:267 down to 251:268.
$ws.Rows.Item(250).Insert()

# Populate the newly inserted row 250 with the weekly Betarraga / Talca record.
$ws.Cells.Item(250, 1).Value = 5
$ws.Cells.Item(250, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(250, 3).Value = "Maule"
$ws.Cells.Item(250, 4).Value = 44610
$ws.Cells.Item(250, 5).Value = 7
$ws.Cells.Item(250, 6).Value = 100114014
$ws.Cells.Item(250, 7).Value = "Betarraga"
$ws.Cells.Item(250, 8).Value = "Sin especificar"
$ws.Cells.Item(250, 9).Value = "Primera"
$ws.Cells.Item(250, 10).Value = 3000
$ws.Cells.Item(250, 11).Value = 800
$ws.Cells.Item(250, 12).Value = 800
$ws.Cells.Item(250, 13).Value = 800
$ws.Cells.Item(250, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(250, 15).Value = "Región del Maule"
$ws.Cells.Item(250, 16).Value = 160
$ws.Cells.Item(250, 17).Value = 5
$ws.Cells.Item(250, 18).Value = "Hortaliza"
